$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q1)
$ws.Range("B2").Value = 0.1521215168462522
$ws.Range("C2").Value = 0.4582899675117284
$ws.Range("D2").Value = 0.4263240185482037
$ws.Range("E2").Value = 0.6529349267332876
$ws.Range("F2").Value = 0.6492376802903603
$ws.Range("G2").Value = 23

# Row 3 (Q2)
$ws.Range("B3").Value = 0.7108428339446126
$ws.Range("C3").Value = 0.8754549345972755
$ws.Range("D3").Value = 4.635275670739081
$ws.Range("E3").Value = 2.152969036177502
$ws.Range("F3").Value = 2.077908663797673
$ws.Range("G3").Value = 23

# Row 4 (Q3)
$ws.Range("B4").Value = 0.2895745887097694
$ws.Range("C4").Value = 1.240496959324251
$ws.Range("D4").Value = 7.595578471819104
$ws.Range("E4").Value = 2.756007705326512
$ws.Range("F4").Value = 2.802350277210764
$ws.Range("G4").Value = 23

# Row 5 (Q4)
$ws.Range("B5").Value = 0.3238621390752919
$ws.Range("C5").Value = 1.24951236653669
$ws.Range("D5").Value = 8.076301587591825
$ws.Range("E5").Value = 2.841883457777927
$ws.Range("F5").Value = 2.886823850443001
$ws.Range("G5").Value = 23

# Row 6 (Q5)
$ws.Range("B6").Value = 0.3759505404624885
$ws.Range("C6").Value = 1.26396980525482
$ws.Range("D6").Value = 7.453626575147394
$ws.Range("E6").Value = 2.730133069128205
$ws.Range("F6").Value = 2.7648986388152
$ws.Range("G6").Value = 23

# Row 7 (Q6)
$ws.Range("B7").Value = 0.341963470112216
$ws.Range("C7").Value = 1.389240385587007
$ws.Range("D7").Value = 7.894718527440571
$ws.Range("E7").Value = 2.809754175624724
$ws.Range("F7").Value = 2.851546062215722
$ws.Range("G7").Value = 23

# Row 8 (Q7)
$ws.Range("B8").Value = 0.2502161759457917
$ws.Range("C8").Value = 1.450552389820832
$ws.Range("D8").Value = 8.661950858152869
$ws.Range("E8").Value = 2.943119239540401
$ws.Range("F8").Value = 2.998369880143227
$ws.Range("G8").Value = 23

# Row 9 (Q8)
$ws.Range("B9").Value = 0.298823989417038
$ws.Range("C9").Value = 1.462549152856695
$ws.Range("D9").Value = 8.060417993749848
$ws.Range("E9").Value = 2.839087528370664
$ws.Range("F9").Value = 2.886770852395387
$ws.Range("G9").Value = 23

# Row 10 (Q9)
$ws.Range("B10").Value = 0.2664666779633617
$ws.Range("C10").Value = 1.439819604875121
$ws.Range("D10").Value = 8.41458104426475
$ws.Range("E10").Value = 2.900789727688781
$ws.Range("F10").Value = 2.953443758313014
$ws.Range("G10").Value = 23

# Row 11
$ws.Range("B11").Value = 0.1601520388162171
$ws.Range("C11").Value = 1.336129065843578
$ws.Range("D11").Value = 8.162906777384929
$ws.Range("E11").Value = 2.857080113924867
$ws.Range("F11").Value = 2.916699070887114
$ws.Range("G11").Value = 23
